# Apply the "cryptos list" refresh committed by the GitHub Actions job.
# Source data is stored as literal (inline/shared) strings, including values
# that look numeric ("0.9987", "1.0000", ...). Excel's Range.Value setter
# auto-converts number-like text to a real number (losing e.g. trailing
# zeros), so for those cells we force the cell to Text format first with
# NumberFormat = "@" before assigning the literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.339.35"
$ws.Range("E2").Value = "  -0.14%  "

# Row 3
$ws.Range("D3").Value = "1.841.38"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.37"
$ws.Range("E5").Value = "  -0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6278"
$ws.Range("E6").Value = "  -0.61%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07433"
$ws.Range("E8").Value = "  -0.66%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2889"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.92"
$ws.Range("E10").Value = "  +1.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07729"
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").Value = "1.831.26"
$ws.Range("E12").Value = "  -0.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.961"
$ws.Range("E13").Value = "  -0.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6752"
$ws.Range("E14").Value = "  -0.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001023"
$ws.Range("E15").Value = "  +0.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.45"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.228"
$ws.Range("E17").Value = "  +1.41%  "

# Row 18
$ws.Range("D18").Value = "29.394.79"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.62"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.351"
$ws.Range("E22").Value = "  -1.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.01"
$ws.Range("E24").Value = "  -0.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.485"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1345"
$ws.Range("E26").Value = "  -2.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.38"
$ws.Range("E27").Value = "  -0.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07144"
$ws.Range("E28").Value = "  +12.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("E29").Value = "  +5.76%  "

# Row 30
$ws.Range("E30").Value = "  +0.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.050"
$ws.Range("E31").Value = "  -1.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.031"
$ws.Range("E32").Value = "  -0.55%  "

# Row 33
$ws.Range("E33").Value = "  +0.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.137"
$ws.Range("E34").Value = "  -0.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6957"
$ws.Range("E35").Value = "  -0.30%  "

# Row 36
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01838"
$ws.Range("E37").Value = "  +1.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.806"
$ws.Range("E38").Value = "  -1.00%  "

# Row 39
$ws.Range("D39").Value = "1.234.14"
$ws.Range("E39").Value = "  -1.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.798"
$ws.Range("E40").Value = "  +3.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9261"
$ws.Range("E41").Value = "  +2.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.0000"
$ws.Range("E42").Value = "  +0.13%  "

# Row 43
$ws.Range("D43").Value = "2.002.22"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.57"
$ws.Range("E44").Value = "  -0.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.24"
$ws.Range("E45").Value = "  -1.70%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.013"
$ws.Range("E46").Value = "  -0.50%  "

# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.709"
$ws.Range("E47").Value = "  +1.31%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.882"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1137"
$ws.Range("E49").Value = "  -3.12%  "

# Row 50
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3907"
$ws.Range("E50").Value = "  -0.69%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05669"
$ws.Range("E51").Value = "  -0.65%  "
